# Update the MitsosBarton2006Ex314 "alpha_zero" stationary-point generator
# worksheet values (leader/follower restriction expressions & evaluations,
# the modified point, and the bf/BF vectors).
#
# NOTE: the workbook has two worksheets whose names differ only by case
# ("Vector_bf" and "Vector_BF"); Worksheets.Item("name") resolves
# case-insensitively, so sheets are addressed by their 1-based tab index
# instead of by name to avoid writing into the wrong sheet.
#   1: Funciones_Objetivo
#   2: Restricciones_del_lider
#   3: Restricciones_del_follower
#   4: Punto_modificado
#   5: Vector_bf
#   6: Vector_BF
#   7: Vector_Alpha

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Range, [string]$Value)
    # Force the cell to stay a text value (matching the workbook's existing
    # shared-string cells) even when $Value looks like a plain number -
    # otherwise Excel auto-converts it to a numeric cell.
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

# --- Restricciones_del_lider ---
$ws = $wb.Worksheets.Item(2)
Set-TextValue $ws.Range("A2") "-2.1 + x"
Set-TextValue $ws.Range("B2") "1.1"
Set-TextValue $ws.Range("D2") "0.74"
Set-TextValue $ws.Range("A3") "2.1 - x"
Set-TextValue $ws.Range("B3") "-3.1"
Set-TextValue $ws.Range("D3") "0.27"

# --- Restricciones_del_follower ---
$ws = $wb.Worksheets.Item(3)
Set-TextValue $ws.Range("A2") "-3.3000000000000003 + y"
Set-TextValue $ws.Range("B2") "2.3000000000000003"
Set-TextValue $ws.Range("D2") "0.22"
Set-TextValue $ws.Range("E2") "0.8"
Set-TextValue $ws.Range("F2") "6.2"
Set-TextValue $ws.Range("A3") "3.3000000000000007 - y"
Set-TextValue $ws.Range("B3") "-4.300000000000001"
Set-TextValue $ws.Range("D3") "0.66"
Set-TextValue $ws.Range("E3") "0"
Set-TextValue $ws.Range("F3") "7.800000000000001"

# --- Punto_modificado ---
$ws = $wb.Worksheets.Item(4)
Set-TextValue $ws.Range("A2") "2.1"
Set-TextValue $ws.Range("B2") "3.3000000000000003"

# --- Vector_bf ---
$ws = $wb.Worksheets.Item(5)
Set-TextValue $ws.Range("A2") "-8.350000000000003"

# --- Vector_BF ---
$ws = $wb.Worksheets.Item(6)
Set-TextValue $ws.Range("A2") "-4.17"
Set-TextValue $ws.Range("A3") "-7.4"
